$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 556
$ws.Range("F6").Value = 1615
$ws.Range("F10").Value = 2721
$ws.Range("F11").Value = 2721
$ws.Range("F13").Value = 1834
$ws.Range("F15").Value = 308
$ws.Range("F16").Value = 714
$ws.Range("F17").Value = 8
$ws.Range("F18").Value = 6298
$ws.Range("F19").Value = 243
$ws.Range("F20").Value = 92
$ws.Range("F22").Value = 3408
$ws.Range("F23").Value = 884
$ws.Range("F27").Value = 2477
$ws.Range("F29").Value = 384
$ws.Range("F33").Value = 1320
$ws.Range("F35").Value = 13
$ws.Range("F39").Value = 1513
$ws.Range("F40").Value = 33
$ws.Range("F41").Value = 1468

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 96
$ws.Range("F10").Value = 220
$ws.Range("F13").Value = 81
$ws.Range("F17").Value = 341
$ws.Range("G18").Value = 280

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 930
$ws.Range("F6").Value = 47

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 930
$ws.Range("F9").Value = 556
$ws.Range("F10").Value = 47
$ws.Range("F16").Value = 96
$ws.Range("F19").Value = 2721
$ws.Range("F25").Value = 308
$ws.Range("F26").Value = 714
$ws.Range("F27").Value = 6298
$ws.Range("F28").Value = 243
$ws.Range("F29").Value = 92
$ws.Range("F34").Value = 2477
$ws.Range("F35").Value = 384
$ws.Range("F37").Value = 1320
$ws.Range("F40").Value = 341
$ws.Range("G41").Value = 280
$ws.Range("F48").Value = 1513
$ws.Range("F49").Value = 33
